# Lisää paikkausten oikeudet rooliexceliin
# Fill in the role-rights matrix for the two new "Paikkaukset" (patching) rows
# (row 53 = "Paikkaukset / Toteumat", row 54 = "Paikkaukset / Kustannukset")
# on the "Oikeudet" sheet, matching the rights already used by the sibling
# rows in the same block.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oikeudet")

# --- Row 53: "Paikkaukset / Toteumat" ---------------------------------
$ws.Range("D53").Value = "R*,W*"
$ws.Range("E53").Value = "R*,W"
$ws.Range("F53").Value = "R*"
$ws.Range("G53").Value = ""
$ws.Range("H53").Value = "R*"
$ws.Range("I53").Value = "R*"
$ws.Range("J53").Value = "R,W"
$ws.Range("K53").Value = "R*,W*"
$ws.Range("L53").Value = "R*,W"
$ws.Range("M53").Value = ""
$ws.Range("N53").Value = "R*"
$ws.Range("O53").Value = "R*"
$ws.Range("P53").Value = ""
$ws.Range("Q53").Value = "R,W"
$ws.Range("R53").Value = "R+,W+"
$ws.Range("S53").Value = "R,W"
$ws.Range("T53").Value = "R+"
$ws.Range("U53").Value = "R"
$ws.Range("V53").Value = "R+"
$ws.Range("W53").Value = ""
$ws.Range("X53").Value = ""

# --- Row 54: "Paikkaukset / Kustannukset" ------------------------------
$ws.Range("D54").Value = "R*,W*"
$ws.Range("E54").Value = "R*,W"
$ws.Range("F54").Value = "R*"
$ws.Range("G54").Value = ""
$ws.Range("H54").Value = "R*"
$ws.Range("I54").Value = "R*"
$ws.Range("J54").Value = "R,W"
$ws.Range("K54").Value = "R*,W*"
$ws.Range("L54").Value = "R*,W"
$ws.Range("M54").Value = ""
$ws.Range("N54").Value = "R*"
$ws.Range("O54").Value = "R*"
$ws.Range("P54").Value = ""
$ws.Range("Q54").Value = "R,W"
$ws.Range("R54").Value = "R+,W+"
$ws.Range("S54").Value = "R,W"
$ws.Range("T54").Value = "R+"
$ws.Range("U54").Value = "R"
$ws.Range("V54").Value = "R+"
$ws.Range("W54").Value = ""
$ws.Range("X54").Value = ""

# --- Borders: give every populated cell a full thin box (matching the
#     rest of the rights matrix), by adding the left edge across the row
#     (and the inside vertical edges between columns) - row 54 also needs
#     its top edge since it previously lacked one.
$row53 = $ws.Range("D53:X53")
$row53.Borders.Item(7).LineStyle = 1
$row53.Borders.Item(7).Weight = 2
$row53.Borders.Item(11).LineStyle = 1
$row53.Borders.Item(11).Weight = 2

$row54 = $ws.Range("D54:X54")
$row54.Borders.Item(7).LineStyle = 1
$row54.Borders.Item(7).Weight = 2
$row54.Borders.Item(11).LineStyle = 1
$row54.Borders.Item(11).Weight = 2
$row54.Borders.Item(8).LineStyle = 1
$row54.Borders.Item(8).Weight = 2
